$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: copy row 568 (values + number formats) down into the 14 new rows
# so dates get style s="1" and nothing pollutes styles.xml.
$srcRow = $ws.Range("A568:T568")
for ($r = 569; $r -le 582; $r++) {
    $destRow = $ws.Range("A" + $r + ":T" + $r)
    $srcRow.Copy($destRow)
}

# Step 2: overwrite per-row cell values/text/formula per the source data.

# Row 569
$ws.Cells.Item(569, 1).Value = 45868
$ws.Cells.Item(569, 2).Value = "Flowering"
$ws.Cells.Item(569, 3).Value = "Large"
$ws.Cells.Item(569, 4).Value = 70
$ws.Cells.Item(569, 5).Value = 91
$ws.Cells.Item(569, 6).Formula = "=ABS(D569-E569)"
$ws.Cells.Item(569, 7).Value = 0
$ws.Cells.Item(569, 8).Value = 0.1
$ws.Cells.Item(569, 9).Value = "No"
$ws.Cells.Item(569, 10).Value = 2
$ws.Cells.Item(569, 11).Value = "Bright"
$ws.Cells.Item(569, 12).Value = 7
$ws.Cells.Item(569, 13).Value = 0.55000000000000004
$ws.Cells.Item(569, 14).Value = 72
$ws.Cells.Item(569, 15).Value = 30.03
$ws.Cells.Item(569, 16).Value = 8
$ws.Cells.Item(569, 17).Value = 0.65
$ws.Cells.Item(569, 18).Value = 8.8000000000000007
$ws.Cells.Item(569, 19).Value = 67
$ws.Cells.Item(569, 20).Value = 17

# Row 570
$ws.Cells.Item(570, 1).Value = 45868
$ws.Cells.Item(570, 2).Value = "Nonflowering"
$ws.Cells.Item(570, 3).Value = "Medium"
$ws.Cells.Item(570, 4).Value = 70
$ws.Cells.Item(570, 5).Value = 91
$ws.Cells.Item(570, 6).Formula = "=ABS(D570-E570)"
$ws.Cells.Item(570, 7).Value = 0
$ws.Cells.Item(570, 8).Value = 0
$ws.Cells.Item(570, 9).Value = "No"
$ws.Cells.Item(570, 10).Value = 3
$ws.Cells.Item(570, 11).Value = "Bright"
$ws.Cells.Item(570, 12).Value = 7
$ws.Cells.Item(570, 13).Value = 0.55000000000000004
$ws.Cells.Item(570, 14).Value = 72
$ws.Cells.Item(570, 15).Value = 30.03
$ws.Cells.Item(570, 16).Value = 8
$ws.Cells.Item(570, 17).Value = 0.65
$ws.Cells.Item(570, 18).Value = 8.8000000000000007
$ws.Cells.Item(570, 19).Value = 67
$ws.Cells.Item(570, 20).Value = 17

# Row 571
$ws.Cells.Item(571, 1).Value = 45868
$ws.Cells.Item(571, 2).Value = "Nonflowering"
$ws.Cells.Item(571, 3).Value = "Small"
$ws.Cells.Item(571, 4).Value = 70
$ws.Cells.Item(571, 5).Value = 91
$ws.Cells.Item(571, 6).Formula = "=ABS(D571-E571)"
$ws.Cells.Item(571, 7).Value = 0
$ws.Cells.Item(571, 8).Value = 0.2
$ws.Cells.Item(571, 9).Value = "No"
$ws.Cells.Item(571, 10).Value = 3
$ws.Cells.Item(571, 11).Value = "Neutral"
$ws.Cells.Item(571, 12).Value = 7
$ws.Cells.Item(571, 13).Value = 0.55000000000000004
$ws.Cells.Item(571, 14).Value = 72
$ws.Cells.Item(571, 15).Value = 30.03
$ws.Cells.Item(571, 16).Value = 8
$ws.Cells.Item(571, 17).Value = 0.65
$ws.Cells.Item(571, 18).Value = 8.8000000000000007
$ws.Cells.Item(571, 19).Value = 67
$ws.Cells.Item(571, 20).Value = 17

# Row 572
$ws.Cells.Item(572, 1).Value = 45868
$ws.Cells.Item(572, 2).Value = "Nonflowering"
$ws.Cells.Item(572, 3).Value = "Medium"
$ws.Cells.Item(572, 4).Value = 70
$ws.Cells.Item(572, 5).Value = 91
$ws.Cells.Item(572, 6).Formula = "=ABS(D572-E572)"
$ws.Cells.Item(572, 7).Value = 0
$ws.Cells.Item(572, 8).Value = 0.3
$ws.Cells.Item(572, 9).Value = "No"
$ws.Cells.Item(572, 10).Value = 3
$ws.Cells.Item(572, 11).Value = "Neutral"
$ws.Cells.Item(572, 12).Value = 7
$ws.Cells.Item(572, 13).Value = 0.55000000000000004
$ws.Cells.Item(572, 14).Value = 72
$ws.Cells.Item(572, 15).Value = 30.03
$ws.Cells.Item(572, 16).Value = 8
$ws.Cells.Item(572, 17).Value = 0.65
$ws.Cells.Item(572, 18).Value = 8.8000000000000007
$ws.Cells.Item(572, 19).Value = 67
$ws.Cells.Item(572, 20).Value = 17

# Row 573
$ws.Cells.Item(573, 1).Value = 45868
$ws.Cells.Item(573, 2).Value = "Nonflowering"
$ws.Cells.Item(573, 3).Value = "Medium"
$ws.Cells.Item(573, 4).Value = 70
$ws.Cells.Item(573, 5).Value = 91
$ws.Cells.Item(573, 6).Formula = "=ABS(D573-E573)"
$ws.Cells.Item(573, 7).Value = 0
$ws.Cells.Item(573, 8).Value = 0.5
$ws.Cells.Item(573, 9).Value = "No"
$ws.Cells.Item(573, 10).Value = 3
$ws.Cells.Item(573, 11).Value = "Bright"
$ws.Cells.Item(573, 12).Value = 7
$ws.Cells.Item(573, 13).Value = 0.55000000000000004
$ws.Cells.Item(573, 14).Value = 72
$ws.Cells.Item(573, 15).Value = 30.03
$ws.Cells.Item(573, 16).Value = 8
$ws.Cells.Item(573, 17).Value = 0.65
$ws.Cells.Item(573, 18).Value = 8.8000000000000007
$ws.Cells.Item(573, 19).Value = 67
$ws.Cells.Item(573, 20).Value = 17

# Row 574
$ws.Cells.Item(574, 1).Value = 45868
$ws.Cells.Item(574, 2).Value = "Nonflowering"
$ws.Cells.Item(574, 3).Value = "Large"
$ws.Cells.Item(574, 4).Value = 70
$ws.Cells.Item(574, 5).Value = 91
$ws.Cells.Item(574, 6).Formula = "=ABS(D574-E574)"
$ws.Cells.Item(574, 7).Value = 0
$ws.Cells.Item(574, 8).Value = 0
$ws.Cells.Item(574, 9).Value = "No"
$ws.Cells.Item(574, 10).Value = 4
$ws.Cells.Item(574, 11).Value = "Neutral"
$ws.Cells.Item(574, 12).Value = 7
$ws.Cells.Item(574, 13).Value = 0.55000000000000004
$ws.Cells.Item(574, 14).Value = 72
$ws.Cells.Item(574, 15).Value = 30.03
$ws.Cells.Item(574, 16).Value = 8
$ws.Cells.Item(574, 17).Value = 0.65
$ws.Cells.Item(574, 18).Value = 8.8000000000000007
$ws.Cells.Item(574, 19).Value = 67
$ws.Cells.Item(574, 20).Value = 17

# Row 575
$ws.Cells.Item(575, 1).Value = 45868
$ws.Cells.Item(575, 2).Value = "Tree"
$ws.Cells.Item(575, 3).Value = "Medium"
$ws.Cells.Item(575, 4).Value = 70
$ws.Cells.Item(575, 5).Value = 91
$ws.Cells.Item(575, 6).Formula = "=ABS(D575-E575)"
$ws.Cells.Item(575, 7).Value = 0
$ws.Cells.Item(575, 8).Value = 1
$ws.Cells.Item(575, 9).Value = "No"
$ws.Cells.Item(575, 10).Value = 1
$ws.Cells.Item(575, 11).Value = "Bright"
$ws.Cells.Item(575, 12).Value = 7
$ws.Cells.Item(575, 13).Value = 0.55000000000000004
$ws.Cells.Item(575, 14).Value = 72
$ws.Cells.Item(575, 15).Value = 30.03
$ws.Cells.Item(575, 16).Value = 8
$ws.Cells.Item(575, 17).Value = 0.65
$ws.Cells.Item(575, 18).Value = 8.8000000000000007
$ws.Cells.Item(575, 19).Value = 67
$ws.Cells.Item(575, 20).Value = 17

# Row 576
$ws.Cells.Item(576, 1).Value = 45869
$ws.Cells.Item(576, 2).Value = "Flowering"
$ws.Cells.Item(576, 3).Value = "Large"
$ws.Cells.Item(576, 4).Value = 62
$ws.Cells.Item(576, 5).Value = 80
$ws.Cells.Item(576, 6).Formula = "=ABS(D576-E576)"
$ws.Cells.Item(576, 7).Value = 0
$ws.Cells.Item(576, 8).Value = 0
$ws.Cells.Item(576, 9).Value = "No"
$ws.Cells.Item(576, 10).Value = 2
$ws.Cells.Item(576, 11).Value = "Neutral"
$ws.Cells.Item(576, 12).Value = 3
$ws.Cells.Item(576, 13).Value = 0.82
$ws.Cells.Item(576, 14).Value = 66
$ws.Cells.Item(576, 15).Value = 30.01
$ws.Cells.Item(576, 16).Value = 16
$ws.Cells.Item(576, 17).Value = 0.89
$ws.Cells.Item(576, 18).Value = 8.1
$ws.Cells.Item(576, 19).Value = 54
$ws.Cells.Item(576, 20).Value = 7

# Row 577
$ws.Cells.Item(577, 1).Value = 45869
$ws.Cells.Item(577, 2).Value = "Nonflowering"
$ws.Cells.Item(577, 3).Value = "Medium"
$ws.Cells.Item(577, 4).Value = 62
$ws.Cells.Item(577, 5).Value = 80
$ws.Cells.Item(577, 6).Formula = "=ABS(D577-E577)"
$ws.Cells.Item(577, 7).Value = 0
$ws.Cells.Item(577, 8).Value = 0.2
$ws.Cells.Item(577, 9).Value = "No"
$ws.Cells.Item(577, 10).Value = 3
$ws.Cells.Item(577, 11).Value = "Dark"
$ws.Cells.Item(577, 12).Value = 3
$ws.Cells.Item(577, 13).Value = 0.82
$ws.Cells.Item(577, 14).Value = 66
$ws.Cells.Item(577, 15).Value = 30.01
$ws.Cells.Item(577, 16).Value = 16
$ws.Cells.Item(577, 17).Value = 0.89
$ws.Cells.Item(577, 18).Value = 8.1
$ws.Cells.Item(577, 19).Value = 54
$ws.Cells.Item(577, 20).Value = 7

# Row 578
$ws.Cells.Item(578, 1).Value = 45869
$ws.Cells.Item(578, 2).Value = "Nonflowering"
$ws.Cells.Item(578, 3).Value = "Small"
$ws.Cells.Item(578, 4).Value = 62
$ws.Cells.Item(578, 5).Value = 80
$ws.Cells.Item(578, 6).Formula = "=ABS(D578-E578)"
$ws.Cells.Item(578, 7).Value = 0
$ws.Cells.Item(578, 8).Value = 0.4
$ws.Cells.Item(578, 9).Value = "No"
$ws.Cells.Item(578, 10).Value = 3
$ws.Cells.Item(578, 11).Value = "Dark"
$ws.Cells.Item(578, 12).Value = 3
$ws.Cells.Item(578, 13).Value = 0.82
$ws.Cells.Item(578, 14).Value = 66
$ws.Cells.Item(578, 15).Value = 30.01
$ws.Cells.Item(578, 16).Value = 16
$ws.Cells.Item(578, 17).Value = 0.89
$ws.Cells.Item(578, 18).Value = 8.1
$ws.Cells.Item(578, 19).Value = 54
$ws.Cells.Item(578, 20).Value = 7

# Row 579
$ws.Cells.Item(579, 1).Value = 45869
$ws.Cells.Item(579, 2).Value = "Nonflowering"
$ws.Cells.Item(579, 3).Value = "Medium"
$ws.Cells.Item(579, 4).Value = 62
$ws.Cells.Item(579, 5).Value = 80
$ws.Cells.Item(579, 6).Formula = "=ABS(D579-E579)"
$ws.Cells.Item(579, 7).Value = 0
$ws.Cells.Item(579, 8).Value = 0.1
$ws.Cells.Item(579, 9).Value = "No"
$ws.Cells.Item(579, 10).Value = 3
$ws.Cells.Item(579, 11).Value = "Bright"
$ws.Cells.Item(579, 12).Value = 3
$ws.Cells.Item(579, 13).Value = 0.82
$ws.Cells.Item(579, 14).Value = 66
$ws.Cells.Item(579, 15).Value = 30.01
$ws.Cells.Item(579, 16).Value = 16
$ws.Cells.Item(579, 17).Value = 0.89
$ws.Cells.Item(579, 18).Value = 8.1
$ws.Cells.Item(579, 19).Value = 54
$ws.Cells.Item(579, 20).Value = 7

# Row 580
$ws.Cells.Item(580, 1).Value = 45869
$ws.Cells.Item(580, 2).Value = "Nonflowering"
$ws.Cells.Item(580, 3).Value = "Medium"
$ws.Cells.Item(580, 4).Value = 62
$ws.Cells.Item(580, 5).Value = 80
$ws.Cells.Item(580, 6).Formula = "=ABS(D580-E580)"
$ws.Cells.Item(580, 7).Value = 0
$ws.Cells.Item(580, 8).Value = 0.5
$ws.Cells.Item(580, 9).Value = "No"
$ws.Cells.Item(580, 10).Value = 3
$ws.Cells.Item(580, 11).Value = "Neutral"
$ws.Cells.Item(580, 12).Value = 3
$ws.Cells.Item(580, 13).Value = 0.82
$ws.Cells.Item(580, 14).Value = 66
$ws.Cells.Item(580, 15).Value = 30.01
$ws.Cells.Item(580, 16).Value = 16
$ws.Cells.Item(580, 17).Value = 0.89
$ws.Cells.Item(580, 18).Value = 8.1
$ws.Cells.Item(580, 19).Value = 54
$ws.Cells.Item(580, 20).Value = 7

# Row 581
$ws.Cells.Item(581, 1).Value = 45869
$ws.Cells.Item(581, 2).Value = "Nonflowering"
$ws.Cells.Item(581, 3).Value = "Large"
$ws.Cells.Item(581, 4).Value = 62
$ws.Cells.Item(581, 5).Value = 80
$ws.Cells.Item(581, 6).Formula = "=ABS(D581-E581)"
$ws.Cells.Item(581, 7).Value = 0
$ws.Cells.Item(581, 8).Value = 0
$ws.Cells.Item(581, 9).Value = "No"
$ws.Cells.Item(581, 10).Value = 4
$ws.Cells.Item(581, 11).Value = "Neutral"
$ws.Cells.Item(581, 12).Value = 3
$ws.Cells.Item(581, 13).Value = 0.82
$ws.Cells.Item(581, 14).Value = 66
$ws.Cells.Item(581, 15).Value = 30.01
$ws.Cells.Item(581, 16).Value = 16
$ws.Cells.Item(581, 17).Value = 0.89
$ws.Cells.Item(581, 18).Value = 8.1
$ws.Cells.Item(581, 19).Value = 54
$ws.Cells.Item(581, 20).Value = 7

# Row 582
$ws.Cells.Item(582, 1).Value = 45869
$ws.Cells.Item(582, 2).Value = "Tree"
$ws.Cells.Item(582, 3).Value = "Medium"
$ws.Cells.Item(582, 4).Value = 62
$ws.Cells.Item(582, 5).Value = 80
$ws.Cells.Item(582, 6).Formula = "=ABS(D582-E582)"
$ws.Cells.Item(582, 7).Value = 0
$ws.Cells.Item(582, 8).Value = 0.5
$ws.Cells.Item(582, 9).Value = "No"
$ws.Cells.Item(582, 10).Value = 1
$ws.Cells.Item(582, 11).Value = "Dark"
$ws.Cells.Item(582, 12).Value = 3
$ws.Cells.Item(582, 13).Value = 0.82
$ws.Cells.Item(582, 14).Value = 66
$ws.Cells.Item(582, 15).Value = 30.01
$ws.Cells.Item(582, 16).Value = 16
$ws.Cells.Item(582, 17).Value = 0.89
$ws.Cells.Item(582, 18).Value = 8.1
$ws.Cells.Item(582, 19).Value = 54
$ws.Cells.Item(582, 20).Value = 7


# Step 3: update the sheet view to match the workbook's saved scroll/selection
# state (selection -> Q569:Q575, matching the post-edit cursor position).
$null = $ws.Range("Q569:Q575").Select()
